# SubscriptionData.xlsx - "Run Tests With Cucumber" edit
#
# Reproduces (as closely as the COM surface allows):
#   1. KSA sheet   (sheet1): selection -> C12:D12 (active cell C12),
#                             B2 value 15 -> 20, orientation -> portrait
#   2. Bahrain sheet (sheet2): selection -> A3, orientation -> portrait
#   3. Kuwait sheet  (sheet3): selection -> A1:C4 (active cell C4),
#                              orientation -> portrait

$wb = $excel.ActiveWorkbook

$wsKSA     = $wb.Worksheets.Item(1)
$wsBahrain = $wb.Worksheets.Item(2)
$wsKuwait  = $wb.Worksheets.Item(3)

# --- Bahrain sheet: just move the selection to A3, set page orientation ---
$wsBahrain.Select()
$wsBahrain.Range("A3").Select()
$wsBahrain.PageSetup.Orientation = 1

# --- Kuwait sheet: change the value in B2, select A1:C4, set orientation ---
$wsKuwait.Select()
$wsKuwait.Range("A1:C4").Select()
$wsKuwait.PageSetup.Orientation = 1

# --- KSA sheet: update B2 value, move selection to C12:D12, set orientation ---
$wsKSA.Select()
$wsKSA.Range("B2").Value = 20
$wsKSA.Range("C12:D12").Select()
$wsKSA.PageSetup.Orientation = 1

# Leave the KSA sheet active/selected, matching the workbook's saved state.
$wsKSA.Select()
$wsKSA.Range("C12:D12").Select()
